$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.237.84'
$ws.Range("E2").Value = '  +1.81%  '

$ws.Range("D3").Value = '2.795.99'
$ws.Range("E3").Value = '  +1.80%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '345.80'
$ws.Range("E5").Value = '  +4.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '116.96'
$ws.Range("E6").Value = '  +1.55%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.553'
$ws.Range("E7").Value = '  +4.24%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.584'
$ws.Range("E9").Value = '  +2.50%  '

$ws.Range("E10").Value = '  +4.24%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0858'
$ws.Range("E11").Value = '  +3.52%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.16'
$ws.Range("E12").Value = '  -0.57%  '

$ws.Range("E13").Value = '  +1.95%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.78'
$ws.Range("E14").Value = '  +1.06%  '

$ws.Range("D15").Value = '3.238.66'
$ws.Range("E15").Value = '  +1.97%  '

$ws.Range("D16").Value = '2.823.29'
$ws.Range("E16").Value = '  +2.98%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.892'
$ws.Range("E17").Value = '  +0.80%  '

$ws.Range("D18").Value = '52.160.94'
$ws.Range("E18").Value = '  +1.69%  '

$ws.Range("E19").Value = '  +6.86%  '

$ws.Range("E20").Value = '  +3.74%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.41'
$ws.Range("E21").Value = '  -1.15%  '

$ws.Range("D22").Value = '0.0₃0983'
$ws.Range("E22").Value = '  +2.24%  '

$ws.Range("E23").Value = '  -0.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '270.11'
$ws.Range("E24").Value = '  -5.18%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.77'
$ws.Range("E25").Value = '  +6.83%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.69'
$ws.Range("E26").Value = '  -0.60%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  +0.00%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.27'
$ws.Range("E28").Value = '  -0.32%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.23'
$ws.Range("E29").Value = '  +0.46%  '

$ws.Range("E30").Value = '  -0.26%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.11'
$ws.Range("E31").Value = '  -0.84%  '

$ws.Range("E32").Value = '  +0.19%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.74'
$ws.Range("E33").Value = '  +2.10%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0824'
$ws.Range("E34").Value = '  +0.01%  '

$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0411'
$ws.Range("E35").Value = '  +16.51%  '

$ws.Range("E36").Value = '  +1.06%  '

$ws.Range("E37").Value = '  -0.17%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.02'
$ws.Range("E38").Value = '  -1.90%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.99'
$ws.Range("E39").Value = '  -0.61%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.23'
$ws.Range("E40").Value = '  +0.34%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.70'
$ws.Range("E41").Value = '  +21.85%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '23.66'
$ws.Range("E42").Value = '  -0.12%  '

$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '127.94'
$ws.Range("E43").Value = '  -1.21%  '

$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.116'
$ws.Range("E44").Value = '  +2.37%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.31'
$ws.Range("E45").Value = '  +0.86%  '

$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '2.080.17'
$ws.Range("E46").Value = '  -1.31%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.35'
$ws.Range("E47").Value = '  -1.75%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.34'
$ws.Range("E48").Value = '  +5.02%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.981'
$ws.Range("E49").Value = '  +19.00%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.55'
$ws.Range("E50").Value = '  +0.80%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.97'
$ws.Range("E51").Value = '  -0.98%  '
